$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F..V (match-result / odds / url columns). A..E (index/country/
# tournament/season/date) stay put for each physical row.
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-Rows($r1, $r2) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# Row pairs whose F:V (match data) content was swapped between the two rows.
$pairs = @(
    @(7, 8),
    @(81, 82),
    @(85, 86),
    @(92, 93),
    @(101, 102),
    @(112, 113),
    @(116, 117),
    @(119, 120),
    @(128, 129)
)

foreach ($pair in $pairs) {
    Swap-Rows $pair[0] $pair[1]
}

# Append two new match rows (131, 132) after the existing last row (130),
# copying the formatting (styles/number formats) from row 130 first.
$ws.Range("A130:V130").Copy()
$ws.Range("A131:V132").PasteSpecial(-4122)

function Set-RowData($r, $data) {
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("E$r").Value = $data[4]
    $ws.Range("F$r").Value = $data[5]
    $ws.Range("G$r").Value = $data[6]
    $ws.Range("H$r").Value = $data[7]
    $ws.Range("I$r").Value = $data[8]
    $ws.Range("J$r").Value = $data[9]
    $ws.Range("K$r").Value = $data[10]
    $ws.Range("L$r").Value = $data[11]
    $ws.Range("M$r").Value = $data[12]
    $ws.Range("N$r").Value = $data[13]
    $ws.Range("O$r").Value = $data[14]
    $ws.Range("P$r").Value = $data[15]
    $ws.Range("Q$r").Value = $data[16]
    $ws.Range("R$r").Value = $data[17]
    $ws.Range("S$r").Value = $data[18]
    $ws.Range("T$r").Value = $data[19]
    $ws.Range("U$r").Value = $data[20]
    $ws.Range("V$r").Value = $data[21]
}

$row131 = @(
    130, "saudi-arabia", "saudi-professional-league", "2023-2024", 45261.66666666666,
    "Al Feiha", 2, "Damac", 4,
    2.01, "25/11/2023 19:13",
    2.5, "01/12/2023 15:51",
    3.4, "25/11/2023 19:13",
    3.44, "01/12/2023 15:51",
    3.56, "25/11/2023 19:13",
    2.83, "01/12/2023 15:59",
    "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-feiha-damac/UoorV1eS/"
)

$row132 = @(
    131, "saudi-arabia", "saudi-professional-league", "2023-2024", 45261.79166666666,
    "Al Hilal", 3, "Al Nassr", 0,
    1.76, "25/11/2023 18:13",
    1.88, "01/12/2023 18:59",
    4.25, "25/11/2023 18:13",
    4.41, "01/12/2023 18:59",
    3.67, "25/11/2023 18:13",
    3.48, "01/12/2023 18:59",
    "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-hilal-al-nassr/nRX9QL3q/"
)

Set-RowData 131 $row131
Set-RowData 132 $row132
